$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
for ($j=1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    Write-Host "  Shape $j : Id=$($sh.Id) $($sh.Name)"
}
